$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the city_Program variable label (A5)
$ws.Range("A5").Value = "Municipal Recycling Program"

# Fix "multi-family" -> "multifamily" in the recycle-rate description (F15)
$ws.Range("F15").Value = "Total percentage of residential (includes single-family and multifamily) recycled materials that were recovered or diverted from a landfill. The numbers for Seattle were pulled direclty from report; Portland and Los Angeles had to be calculated by hand from available data."

# Add "Bureau" to the Census source citations (F16, F17)
$ws.Range("F16").Value = "Population of city measured in individual residents (US Census Bureau data)."
$ws.Range("F17").Value = "Area of each city measured in square miles (US Census Bureau data)."

# Row 15 auto-fits shorter now that the description text above is shorter
$ws.Rows.Item(15).RowHeight = 75.6

# Move the active selection to D15 (matches the editor's last cursor position)
$ws.Range("D15").Select()
